# Add two new columns: I ("I0") and J ("IF") to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — same style (s=1, bold + border + centered) as the
# existing headers in B1:H1. Copy formatting via PasteSpecial (xlPasteFormats)
# since Range.Style refers to named cell styles, not direct formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row values for I ("I0") and J ("IF"). For every data row (2-33) I0 is 1
# and IF equals the existing H ("IP") value, except row 17 where I0 is 2 and
# IF is 8 (one more than H17's 7).
$iVals = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 1
}
$jVals = @{
    2  = 7
    3  = 5
    4  = 5
    5  = 6
    6  = 6
    7  = 6
    8  = 5
    9  = 5
    10 = 5
    11 = 7
    12 = 3
    13 = 6
    14 = 7
    15 = 7
    16 = 5
    17 = 8
    18 = 5
    19 = 7
    20 = 4
    21 = 7
    22 = 6
    23 = 4
    24 = 6
    25 = 6
    26 = 4
    27 = 5
    28 = 6
    29 = 6
    30 = 5
    31 = 4
    32 = 3
    33 = 2
}

for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 9).Value = $iVals[$row]
    $ws.Cells.Item($row, 10).Value = $jVals[$row]
}
